# Apply updated dSF (column F) values for specific rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7  = -2
    11 = 0
    23 = 4
    36 = 3
    37 = -3
    40 = -4
    42 = -8
    46 = 6
    51 = 7
    52 = -8
    55 = -4
    57 = 9
    61 = -6
    63 = -5
    65 = -5
    67 = 5
    68 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
